$d = $word.ActiveDocument

# 1) "Why " + "Logistic regression" + "?" runs collapse to one run with the
#    same visible text "Why Logistic regression?" -- a no-visible-change
#    replace forces the engine to rebuild/merge the backing runs.
$d.Content.Find.Execute("Why Logistic regression?", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Why Logistic regression?", 2) | Out-Null

# 2) Fix word order typo: "3 the days" -> "the 3 days"
$d.Content.Find.Execute("at the start of 3 the days. It will be coded", $true, $false, $false, $false, $false, `
    $true, 1, $false, "at the start of the 3 days. It will be coded", 2) | Out-Null

# 3) "higher high and" -> "higher highs and"
$d.Content.Find.Execute("if the stock creates 3 consecutive higher high and", $true, $false, $false, $false, $false, `
    $true, 1, $false, "if the stock creates 3 consecutive higher highs and", 2) | Out-Null

# 4) "3 lower lows" -> "3 consecutive lower lows"
$d.Content.Find.Execute("if the stock creates 3 lower lows", $true, $false, $false, $false, $false, `
    $true, 1, $false, "if the stock creates 3 consecutive lower lows", 2) | Out-Null

# 5) ", " + "Bureau of" + " Labor Statistics" runs collapse to one run with
#    the same visible text ", Bureau of Labor Statistics".
$d.Content.Find.Execute(", Bureau of Labor Statistics", $true, $false, $false, $false, $false, `
    $true, 1, $false, ", Bureau of Labor Statistics", 2) | Out-Null

# 6) Remove the stale <w:lastRenderedPageBreak/> marker before "Conclusion".
#    Re-assigning the paragraph's own text rebuilds its run (preserving the
#    bold formatting) and drops the rendering-hint element.
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -eq "Conclusion`r") {
        $p.Range.Text = "Conclusion"
        break
    }
}
